$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the GW27 results for row 24 (week 26)
$ws.Range("B24").Value = 83
$ws.Range("C24").Value = 72
$ws.Range("D24").Value = 82

# Recalculate so the SUM totals in row 41 pick up the new values
$excel.Calculate()

# Update the selection to match the author's final cursor position
$ws.Range("M10").Select()
